$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clean up the "Objet : " prefix from the objet column texts, and
# consolidate the duplicated rows (one per tender with "Date limite de
# remise des plis" row + one with "N/A" row) into a single row per
# tender where column B is always "N/A".

$tenders = @(
    "La réalisation des travaux d’entretien de la chaussée des autoroutes Kenitra – Sidi El Yamani et Khémisset – Meknes « PROGRAMME 2025-2026 »",
    "La réalisation des travaux d’entretien des Bâtiments des axes autoroutiers Nord, Centre et Est",
    "La réalisation des travaux de construction des viaducs sur Oued Cherrat et sur Oued Yquem de l'autoroute Rabat Casablanca Continentale",
    "Mise en conformité de la ligne 60 kV N°107-1 MOHAMMEDIA-ONCF EL MANSOURIA nécessitée par le projet d’aménagement des Routes Provinciales N°3304 et N°3308 pour la desserte du futur Grand Stade HASSAN II - Province",
    "La réalisation des prestations de transport du personnel d’ADM"
)

for ($i = 0; $i -lt $tenders.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $tenders[$i]
    $ws.Cells.Item($row, 2).Value = "N/A"
}

# Remove the now-unused rows 7-11 (previously duplicated second entries
# per tender).
$ws.Range("A7:B11").EntireRow.Delete()
